# Add team record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): AD1 = "Wins", AE1 = "Losses", AF1 = "Ties".
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, centered,
# bordered) by copying the format from the existing last header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-41: every player row carries the team's season record.
$lastRow = 41
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 68  # AD = Wins
    $ws.Cells.Item($r, 31).Value = 94  # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF = Ties
}
